$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6602
    6  = 2013
    7  = 1526
    10 = 409
    11 = 10
    12 = 5631
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
